$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect to make the edits, then re-protect with the
# same settings it had before.
$ws.Unprotect("D382")

# Update the confidential disclosure date text in cell A59 (shared string)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-56
$ws.Range("D2").Value = 0.01545705658527453
$ws.Range("E2").Value = -0.006769930675909919
$ws.Range("D3").Value = 0.05297938018146897
$ws.Range("E3").Value = 0.01201780285185072
$ws.Range("D4").Value = 0.01519946377983439
$ws.Range("E4").Value = -0.009016360001632084
$ws.Range("D5").Value = 0.009449203585211393
$ws.Range("E5").Value = 0.0002625016406352376
$ws.Range("D6").Value = 0.01506322608278288
$ws.Range("E6").Value = -0.007418293792427533
$ws.Range("D7").Value = 0.0191162510619758
$ws.Range("E7").Value = 0.001677084921481997
$ws.Range("D8").Value = 0.003879766842842473
$ws.Range("E8").Value = 0.05750727232042974
$ws.Range("D9").Value = 0.006296451198550954
$ws.Range("E9").Value = 0.03013650062045747
$ws.Range("D10").Value = 0.01419129241949598
$ws.Range("E10").Value = 0.0007668711656441118
$ws.Range("D11").Value = 0.008027296260044723
$ws.Range("E11").Value = 0.001467748165314697
$ws.Range("D12").Value = 0.01460354012545795
$ws.Range("E12").Value = 0.01401273885350318
$ws.Range("D13").Value = 0.003017909932296072
$ws.Range("E13").Value = -0.0003595828838546256
$ws.Range("D14").Value = 0.006122014857264699
$ws.Range("E14").Value = 0.01838440111420625
$ws.Range("D15").Value = 0.01359974052553927
$ws.Range("E15").Value = 0.006479338842975135
$ws.Range("D16").Value = 0.0101838763734762
$ws.Range("E16").Value = -0.001096040553500433
$ws.Range("D17").Value = 0.02202809154465575
$ws.Range("E17").Value = -0.007030656194577989
$ws.Range("D18").Value = 0.008589486046785283
$ws.Range("E18").Value = 0.001061249241965045
$ws.Range("D19").Value = 0.0161902721200661
$ws.Range("E19").Value = -0.007353842382645026
$ws.Range("D20").Value = 0.01129898533443607
$ws.Range("E20").Value = -0.009110367158772892
$ws.Range("D21").Value = 0.007059865986170967
$ws.Range("E21").Value = 0.006174845628859282
$ws.Range("D22").Value = 0.01304136440214357
$ws.Range("E22").Value = 0.003861003861004075
$ws.Range("D23").Value = 0.01920077176342378
$ws.Range("E23").Value = 0.002790373212417041
$ws.Range("D24").Value = 0.009433049775433551
$ws.Range("E24").Value = 0.01382137069869427
$ws.Range("D25").Value = 0.02066943522130635
$ws.Range("E25").Value = -0.001326053042121522
$ws.Range("D26").Value = 0.01264006159992743
$ws.Range("E26").Value = 0.02184103435841966
$ws.Range("D27").Value = 0.02199423365544768
$ws.Range("E27").Value = -0.01059111772732924
$ws.Range("D28").Value = 0.0583354068034907
$ws.Range("E28").Value = -0.006027234169208917
$ws.Range("D29").Value = 0.02083875567273974
$ws.Range("E29").Value = -0.003370029206919756
$ws.Range("D30").Value = 0.03125095576061906
$ws.Range("E30").Value = -0.01599333280419069
$ws.Range("D31").Value = 0.01583250088968891
$ws.Range("E31").Value = -0.0104340265198174
$ws.Range("D32").Value = 0.01368485032945515
$ws.Range("E32").Value = -0.03228583727938006
$ws.Range("D33").Value = 0.01997094572672818
$ws.Range("E33").Value = -0.03834115805946792
$ws.Range("D34").Value = 0.04261964121862633
$ws.Range("E34").Value = 0.02970781063125805
$ws.Range("D35").Value = 0.0107489496617702
$ws.Range("E35").Value = 0.003807545863620643
$ws.Range("D36").Value = 0.009674674803454133
$ws.Range("E36").Value = 0.00439698492462326
$ws.Range("D37").Value = 0.01077592435372168
$ws.Range("E37").Value = 0.002157962883038378
$ws.Range("D38").Value = 0.007133410778480067
$ws.Range("E38").Value = 0.008345286217238401
$ws.Range("D39").Value = 0.01149940419509628
$ws.Range("E39").Value = 0.003761283851554609
$ws.Range("D40").Value = 0.01688311863310959
$ws.Range("E40").Value = -0.005061319836480394
$ws.Range("D41").Value = 0.01703566516690418
$ws.Range("E41").Value = -0.0004149664568779965
$ws.Range("D42").Value = 0.03451551359460945
$ws.Range("E42").Value = -0.003045247527420636
$ws.Range("D43").Value = 0.0112152267455492
$ws.Range("E43").Value = -0.001146139631533427
$ws.Range("D44").Value = 0.02209819473833656
$ws.Range("E44").Value = 0.01539732938976113
$ws.Range("D45").Value = 0.01301342854300762
$ws.Range("E45").Value = -0.01533183231478508
$ws.Range("D46").Value = 0.008022924499625
$ws.Range("E46").Value = -0.0004173767868944633
$ws.Range("D47").Value = 0.01285520802228529
$ws.Range("E47").Value = 0.008509160371624747
$ws.Range("D48").Value = 0.009692471899063502
$ws.Range("E48").Value = -0.00362756952841603
$ws.Range("D49").Value = 0.01526414102973185
$ws.Range("E49").Value = -0.005155330221447429
$ws.Range("D50").Value = 0.008579471304830883
$ws.Range("E50").Value = -0.0132702586471708
$ws.Range("D51").Value = 0.01185038524439603
$ws.Range("E51").Value = -0.007655597534301006
$ws.Range("D52").Value = 0.008806492792726005
$ws.Range("E52").Value = -0.002710971689710018
$ws.Range("D53").Value = 0.009740654279859599
$ws.Range("E53").Value = -0.02268270944741524
$ws.Range("D54").Value = 0.1347068987859621
$ws.Range("E54").Value = 0.00009851246182646101
$ws.Range("D55").Value = 0.04402269726482001
$ws.Range("E55").Value = -0.001795977011494143
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = -0.0003392169830694902

# Re-protect the sheet with its original password
$ws.Protect("D382")
